# Kết quả test với đầu vào có và không có Volumes
# Adds two "ANN.NET" result tables (Volume / NoVolume) to the ANN sheet
# and removes the extra "period = 30" column (F) from the K-SVMeans sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ANN": add the ANN.NET Volume table (G1:K8) next to the existing
# one, and a second ANN.NET NoVolume table below (rows 10-17).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ANN")

# --- ANN.NET Volume header (G1:K1, merged & centered) ---
$ws.Range("G1").Value = "ANN.NET Volume"
$ws.Range("G1:K1").HorizontalAlignment = -4108
$ws.Range("G1:K1").Merge()

# Column headers (row 2)
$ws.Range("H2").Value = "period = 1"
$ws.Range("I2").Value = "period = 5"
$ws.Range("J2").Value = "period = 10"
$ws.Range("K2").Value = "period = 30"

# Row labels (col G)
$ws.Range("G3").Value = "BT6"
$ws.Range("G4").Value = "DHG"
$ws.Range("G5").Value = "FPT"
$ws.Range("G6").Value = "VIS"
$ws.Range("G7").Value = "VNM"
$ws.Range("G8").Value = "Total"

# Data values
$ws.Range("H3").Value = 68.099999999999994
$ws.Range("I3").Value = 67.78
$ws.Range("J3").Value = 59.33
$ws.Range("K3").Value = 41.3

$ws.Range("H4").Value = 56.98
$ws.Range("I4").Value = 56.75
$ws.Range("J4").Value = 56.52
$ws.Range("K4").Value = 47.22

$ws.Range("H5").Value = 67.02
$ws.Range("I5").Value = 32.619999999999997
$ws.Range("J5").Value = 43.01
$ws.Range("K5").Value = 47.8

$ws.Range("H6").Value = 42.16
$ws.Range("I6").Value = 41.84
$ws.Range("J6").Value = 51.91
$ws.Range("K6").Value = 55.86

$ws.Range("H7").Value = 59.05
$ws.Range("I7").Value = 56.7
$ws.Range("J7").Value = 54.34
$ws.Range("K7").Value = 41.15

# Totals row (averages)
$ws.Range("H8").Formula = "=AVERAGE(H3:H7)"
$ws.Range("I8:K8").Formula = "=AVERAGE(I3:I7)"
$ws.Range("H8").Font.Color = 255
$ws.Range("I8").Font.Color = 255

# --- ANN.NET NoVolume table (rows 10-17, columns A-E) ---
$ws.Range("A10").Value = "ANN.NET NoVolume"
$ws.Range("A10:E10").HorizontalAlignment = -4108
$ws.Range("A10:E10").Merge()

$ws.Range("B11").Value = "period = 1"
$ws.Range("C11").Value = "period = 5"
$ws.Range("D11").Value = "period = 10"
$ws.Range("E11").Value = "period = 30"

$ws.Range("A12").Value = "BT6"
$ws.Range("B12").Value = 72.849999999999994
$ws.Range("C12").Value = 66.099999999999994
$ws.Range("D12").Value = 44.01
$ws.Range("E12").Value = 58.45

$ws.Range("A13").Value = "DHG"
$ws.Range("B13").Value = 56.98
$ws.Range("C13").Value = 56.75
$ws.Range("D13").Value = 44.02
$ws.Range("E13").Value = 31.66

$ws.Range("A14").Value = "FPT"
$ws.Range("B14").Value = 43.61
$ws.Range("C14").Value = 60.96
$ws.Range("D14").Value = 67.2
$ws.Range("E14").Value = 32.409999999999997

$ws.Range("A15").Value = "VIS"
$ws.Range("B15").Value = 49.72
$ws.Range("C15").Value = 38.58
$ws.Range("D15").Value = 60.65
$ws.Range("E15").Value = 52.51

$ws.Range("A16").Value = "VNM"
$ws.Range("B16").Value = 56.46
$ws.Range("C16").Value = 57.14
$ws.Range("D16").Value = 50
$ws.Range("E16").Value = 37.61

$ws.Range("A17").Value = "Total"
$ws.Range("B17").Formula = "=AVERAGE(B12:B16)"
$ws.Range("C17:E17").Formula = "=AVERAGE(C12:C16)"
$ws.Range("B17").Font.Color = 255
$ws.Range("C17").Font.Color = 255

# Fix up the column K width (widened to fit the new "period = 30" header)
$ws.Columns.Item(11).ColumnWidth = 11.28515625

# ---------------------------------------------------------------------
# Sheet "K-SVMeans": drop the extra "period = 30" column (F) from the
# second (No Volume) mini table.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("K-SVMeans")
$ws3.Range("F4:F9").Clear()

# ---------------------------------------------------------------------
# Selections / active sheet to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("F16").Select()

$ws2 = $wb.Worksheets.Item("SVM")
$ws2.Range("B9").Select()

$ws3.Range("B8").Select()

$ws.Activate()
$ws.Range("F16").Select()
